$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-5 with new simulation values
# Row 2
$ws.Cells.Item(2, 1).Value = 45137.50694444445
$ws.Cells.Item(2, 2).Value = 14.835
$ws.Cells.Item(2, 3).Value = 9.791
$ws.Cells.Item(2, 4).Value = 3.698
$ws.Cells.Item(2, 5).Value = 32.243
$ws.Cells.Item(2, 6).Value = 24.166
$ws.Cells.Item(2, 7).Value = 11.51
$ws.Cells.Item(2, 8).Value = 34.958
$ws.Cells.Item(2, 9).Value = 18.033
$ws.Cells.Item(2, 10).Value = 7.29
$ws.Cells.Item(2, 11).Value = 10.735
$ws.Cells.Item(2, 12).Value = 12.533
$ws.Cells.Item(2, 13).Value = 13.25
$ws.Cells.Item(2, 14).Value = 3.739
$ws.Cells.Item(2, 15).Value = 11.655
$ws.Cells.Item(2, 16).Value = 16.06
$ws.Cells.Item(2, 17).Value = 10.282
$ws.Cells.Item(2, 18).Value = 3.096
$ws.Cells.Item(2, 19).Value = 1.74
$ws.Cells.Item(2, 20).Value = 170.025
$ws.Cells.Item(2, 21).Value = 32.298
$ws.Cells.Item(2, 22).Value = 10.758
$ws.Cells.Item(2, 23).Value = 20.812
$ws.Cells.Item(2, 24).Value = 10.713
$ws.Cells.Item(2, 25).Value = 2.837
$ws.Cells.Item(2, 26).Value = 18.288
$ws.Cells.Item(2, 27).Value = 9.502
$ws.Cells.Item(2, 28).Value = 8.642
$ws.Cells.Item(2, 29).Value = 10.303
$ws.Cells.Item(2, 30).Value = 12.679
$ws.Cells.Item(2, 31).Value = 3.311
$ws.Cells.Item(2, 32).Value = 31.418
$ws.Cells.Item(2, 33).Value = 5.68
$ws.Cells.Item(2, 34).Value = 13.449

# Row 3
$ws.Cells.Item(3, 1).Value = 45137.51388888889
$ws.Cells.Item(3, 2).Value = 17.722
$ws.Cells.Item(3, 3).Value = 12.713
$ws.Cells.Item(3, 4).Value = 1.821
$ws.Cells.Item(3, 5).Value = 38.836
$ws.Cells.Item(3, 6).Value = 30.739
$ws.Cells.Item(3, 7).Value = 13.829
$ws.Cells.Item(3, 8).Value = 52.446
$ws.Cells.Item(3, 9).Value = 21.524
$ws.Cells.Item(3, 10).Value = 9.318
$ws.Cells.Item(3, 11).Value = 13.64
$ws.Cells.Item(3, 12).Value = 15.408
$ws.Cells.Item(3, 13).Value = 16.338
$ws.Cells.Item(3, 14).Value = 4.468
$ws.Cells.Item(3, 15).Value = 13.911
$ws.Cells.Item(3, 16).Value = 19.603
$ws.Cells.Item(3, 17).Value = 12.011
$ws.Cells.Item(3, 18).Value = 1.475
$ws.Cells.Item(3, 19).Value = 1.063
$ws.Cells.Item(3, 20).Value = 204.414
$ws.Cells.Item(3, 21).Value = 38.895
$ws.Cells.Item(3, 22).Value = 12.84
$ws.Cells.Item(3, 23).Value = 25.763
$ws.Cells.Item(3, 24).Value = 13.462
$ws.Cells.Item(3, 25).Value = 2.449
$ws.Cells.Item(3, 26).Value = 25.919
$ws.Cells.Item(3, 27).Value = 11.342
$ws.Cells.Item(3, 28).Value = 10.191
$ws.Cells.Item(3, 29).Value = 12.018
$ws.Cells.Item(3, 30).Value = 15.909
$ws.Cells.Item(3, 31).Value = 1.246
$ws.Cells.Item(3, 32).Value = 48.023
$ws.Cells.Item(3, 33).Value = 7.053
$ws.Cells.Item(3, 34).Value = 16.053

# Row 4
$ws.Cells.Item(4, 1).Value = 45137.52083333334
$ws.Cells.Item(4, 2).Value = 6.196
$ws.Cells.Item(4, 3).Value = 4.286
$ws.Cells.Item(4, 4).Value = 0.967
$ws.Cells.Item(4, 5).Value = 13.738
$ws.Cells.Item(4, 6).Value = 10.405
$ws.Cells.Item(4, 7).Value = 4.787
$ws.Cells.Item(4, 8).Value = 23.264
$ws.Cells.Item(4, 9).Value = 7.562
$ws.Cells.Item(4, 10).Value = 3.206
$ws.Cells.Item(4, 11).Value = 4.504
$ws.Cells.Item(4, 12).Value = 5.409
$ws.Cells.Item(4, 13).Value = 5.773
$ws.Cells.Item(4, 14).Value = 1.576
$ws.Cells.Item(4, 15).Value = 4.888
$ws.Cells.Item(4, 16).Value = 6.846
$ws.Cells.Item(4, 17).Value = 4.413
$ws.Cells.Item(4, 18).Value = 0.918
$ws.Cells.Item(4, 19).Value = 0.511
$ws.Cells.Item(4, 20).Value = 67.1
$ws.Cells.Item(4, 21).Value = 13.881
$ws.Cells.Item(4, 22).Value = 4.511
$ws.Cells.Item(4, 23).Value = 9.027
$ws.Cells.Item(4, 24).Value = 4.653
$ws.Cells.Item(4, 25).Value = 1.06
$ws.Cells.Item(4, 26).Value = 10.849
$ws.Cells.Item(4, 27).Value = 3.985
$ws.Cells.Item(4, 28).Value = 3.678
$ws.Cells.Item(4, 29).Value = 4.327
$ws.Cells.Item(4, 30).Value = 5.479
$ws.Cells.Item(4, 31).Value = 0.766
$ws.Cells.Item(4, 32).Value = 21.649
$ws.Cells.Item(4, 33).Value = 2.381
$ws.Cells.Item(4, 34).Value = 5.642

# Row 5
$ws.Cells.Item(5, 1).Value = 45137.52777777778
$ws.Cells.Item(5, 2).Value = 11
$ws.Cells.Item(5, 3).Value = 8.01
$ws.Cells.Item(5, 4).Value = 0.92
$ws.Cells.Item(5, 5).Value = 24.16
$ws.Cells.Item(5, 6).Value = 19.25
$ws.Cells.Item(5, 7).Value = 8.59
$ws.Cells.Item(5, 8).Value = 32.63
$ws.Cells.Item(5, 9).Value = 13.38
$ws.Cells.Item(5, 10).Value = 5.83
$ws.Cells.Item(5, 11).Value = 8.54
$ws.Cells.Item(5, 12).Value = 9.62
$ws.Cells.Item(5, 13).Value = 10.24
$ws.Cells.Item(5, 14).Value = 2.78
$ws.Cells.Item(5, 15).Value = 8.65
$ws.Cells.Item(5, 16).Value = 12.2
$ws.Cells.Item(5, 17).Value = 7.45
$ws.Cells.Item(5, 18).Value = 0.74
$ws.Cells.Item(5, 19).Value = 0.56
$ws.Cells.Item(5, 20).Value = 124.27
$ws.Cells.Item(5, 21).Value = 24.12
$ws.Cells.Item(5, 22).Value = 7.98
$ws.Cells.Item(5, 23).Value = 16.02
$ws.Cells.Item(5, 24).Value = 8.42
$ws.Cells.Item(5, 25).Value = 1.45
$ws.Cells.Item(5, 26).Value = 15.84
$ws.Cells.Item(5, 27).Value = 7.05
$ws.Cells.Item(5, 28).Value = 6.33
$ws.Cells.Item(5, 29).Value = 7.44
$ws.Cells.Item(5, 30).Value = 9.96
$ws.Cells.Item(5, 31).Value = 0.55
$ws.Cells.Item(5, 32).Value = 29.57
$ws.Cells.Item(5, 33).Value = 4.4
$ws.Cells.Item(5, 34).Value = 9.98

# Delete old row 6 (data no longer present; dimension shrinks to AH5)
$ws.Rows.Item(6).Delete()

# Update column widths (stored width = ColumnWidth + 0.8333333333333334)
$ws.Range("B1:C1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("G1:G1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("I1:I1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("K1:M1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("O1:Q1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("T1:T1").EntireColumn.ColumnWidth = 8.166666666666666
$ws.Range("V1:V1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("X1:X1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("Z1:AD1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AE1:AE1").EntireColumn.ColumnWidth = 6.166666666666667
$ws.Range("AH1:AH1").EntireColumn.ColumnWidth = 7.166666666666667
